# Sprint #1 Burndown Chart - backlog number tweak
#
# The "Actual Work" series' last data point (C16) is revised: instead of
# dropping 3 points from the prior day (C15-3 = 64), the team logs a
# bigger finish (C15-5.5 = 61.5). Update the formula in place; Excel
# recalculates the cached <v> automatically, and the chart (which is
# sourced straight from Sheet1!$C$2:$C$16) picks up the new value too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C16").Formula = "=C15-5.5"
